$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 16 — this shifts the existing row 16 (and everything
# below it) down by one, turning the old A1:R63 range into A1:R64.
$ws.Rows("16:16").Insert()

# Populate the newly inserted row 16 with this week's new price report.
$ws.Cells.Item(16, 1).Value = 7
$ws.Cells.Item(16, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(16, 3).Value = "Ñuble"
$ws.Cells.Item(16, 4).Value = 44608
$ws.Cells.Item(16, 5).Value = 16
$ws.Cells.Item(16, 6).Value = 100112022
$ws.Cells.Item(16, 7).Value = "Arveja Verde"
$ws.Cells.Item(16, 8).Value = "Sin especificar"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 60
$ws.Cells.Item(16, 11).Value = 22000
$ws.Cells.Item(16, 12).Value = 23000
$ws.Cells.Item(16, 13).Value = 22500
$ws.Cells.Item(16, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(16, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(16, 16).Value = 900
$ws.Cells.Item(16, 17).Value = 25
$ws.Cells.Item(16, 18).Value = "Hortaliza"

# Make sure the new row's date cell keeps the same date format as the rest
# of column D.
$ws.Cells.Item(16, 4).NumberFormat = $ws.Cells.Item(17, 4).NumberFormat
